$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers (e.g. "1.003") must be forced
# to Text format first, otherwise Excel auto-converts the string to a numeric value.
$textCells = @("D4", "D5", "D7", "D8", "D9", "D10", "D11", "D12", "D15", "D16", "D17", "D19", "D20", "D23", "D24", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '29.418.94'
$ws.Range("E2").Value = '  +1.20%  '
$ws.Range("D3").Value = '1.947.02'
$ws.Range("E3").Value = '  +2.61%  '
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("D5").Value = '325.71'
$ws.Range("E5").Value = '  +0.02%  '
$ws.Range("E6").Value = '  +0.25%  '
$ws.Range("D7").Value = '0.4634'
$ws.Range("E7").Value = '  +0.64%  '
$ws.Range("D8").Value = '0.3872'
$ws.Range("E8").Value = '  -0.46%  '
$ws.Range("D9").Value = '46.17'
$ws.Range("E9").Value = '  +0.64%  '
$ws.Range("D10").Value = '0.07835'
$ws.Range("E10").Value = '  -0.22%  '
$ws.Range("D11").Value = '0.9792'
$ws.Range("E11").Value = '  -1.08%  '
$ws.Range("D12").Value = '22.64'
$ws.Range("E12").Value = '  +3.67%  '
$ws.Range("D13").Value = '1.949.62'
$ws.Range("E13").Value = '  +4.21%  '
$ws.Range("E14").Value = '  +0.45%  '
$ws.Range("D15").Value = '5.752'
$ws.Range("E15").Value = '  -0.31%  '
$ws.Range("D16").Value = '0.07031'
$ws.Range("E16").Value = '  +0.39%  '
$ws.Range("D17").Value = '86.79'
$ws.Range("E17").Value = '  -1.20%  '
$ws.Range("E18").Value = '  +0.18%  '
$ws.Range("D19").Value = '0.000009826'
$ws.Range("E19").Value = '  -1.06%  '
$ws.Range("D20").Value = '17.05'
$ws.Range("E20").Value = '  +0.24%  '
$ws.Range("D22").Value = '29.424.77'
$ws.Range("E22").Value = '  +1.21%  '
$ws.Range("D23").Value = '5.468'
$ws.Range("E23").Value = '  +2.92%  '
$ws.Range("D24").Value = '11.08'
$ws.Range("E24").Value = '  -0.04%  '
$ws.Range("D25").Value = '2.176.15'
$ws.Range("E25").Value = '  +3.53%  '
$ws.Range("D26").Value = '2.099'
$ws.Range("E26").Value = '  +0.08%  '
$ws.Range("D27").Value = '157.27'
$ws.Range("E27").Value = '  +0.86%  '
$ws.Range("D28").Value = '19.33'
$ws.Range("E28").Value = '  -0.16%  '
$ws.Range("D29").Value = '5.771'
$ws.Range("E29").Value = '  -1.83%  '
$ws.Range("D30").Value = '118.57'
$ws.Range("E30").Value = '  +0.01%  '
$ws.Range("D31").Value = '1.869'
$ws.Range("E31").Value = '  -0.31%  '
$ws.Range("D32").Value = '0.09377'
$ws.Range("E32").Value = '  +0.55%  '
$ws.Range("D33").Value = '0.8650'
$ws.Range("E33").Value = '  -3.55%  '
$ws.Range("D34").Value = '5.183'
$ws.Range("E34").Value = '  -0.92%  '
$ws.Range("D35").Value = '1.304'
$ws.Range("E35").Value = '  -1.23%  '
$ws.Range("D36").Value = '3.126'
$ws.Range("E36").Value = '  -0.71%  '
$ws.Range("D37").Value = '0.05764'
$ws.Range("E37").Value = '  -0.46%  '
$ws.Range("D38").Value = '1.157'
$ws.Range("E38").Value = '  -1.14%  '
$ws.Range("D39").Value = '0.02085'
$ws.Range("E39").Value = '  +0.16%  '
$ws.Range("D40").Value = '7.700'
$ws.Range("E40").Value = '  +0.31%  '
$ws.Range("D41").Value = '0.5660'
$ws.Range("E41").Value = '  -0.19%  '
$ws.Range("D42").Value = '0.1785'
$ws.Range("E42").Value = '  -0.50%  '
$ws.Range("D43").Value = '9.440'
$ws.Range("E43").Value = '  -2.75%  '
$ws.Range("D44").Value = '0.000002853'
$ws.Range("E44").Value = '  +35.74%  '
$ws.Range("D45").Value = '2.728'
$ws.Range("E45").Value = '  +6.85%  '
$ws.Range("D46").Value = '0.5292'
$ws.Range("E46").Value = '  -0.91%  '
$ws.Range("D47").Value = '11.60'
$ws.Range("E47").Value = '  -2.36%  '
$ws.Range("D48").Value = '2.101'
$ws.Range("E48").Value = '  -5.86%  '
$ws.Range("D49").Value = '0.06875'
$ws.Range("E49").Value = '  -1.84%  '
$ws.Range("D50").Value = '1.815'
$ws.Range("E50").Value = '  -1.71%  '
$ws.Range("D51").Value = '111.49'
$ws.Range("E51").Value = '  -1.07%  '
